# Update the "Forecast Comparison" sheet:
#  - insert a new "Week_Start_Date" column between "Week" and "ASIN"
#  - fill it with the weekly start dates (stored as text, matching source data)
#  - normalize the "Week" labels from zero-padded (W01) to unpadded (W1) for weeks 1-9

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B (ASIN), shifting ASIN..is_holiday_week right by one.
$ws.Columns("B").Insert()

# Header for the newly inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Make sure the date strings are preserved as plain text (not auto-converted to date serials).
$ws.Range("B2:B17").NumberFormat = "@"

$weekStartDates = @(
    "2024-12-15",
    "2024-12-22",
    "2024-12-29",
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $weekStartDates[$i]
}

# Drop the leading zero from the single-digit week labels (W01 -> W1, ... W09 -> W9).
for ($week = 1; $week -le 9; $week++) {
    $row = $week + 1
    $ws.Range("A$row").Value = "W$week"
}
